$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sim")

# District filename and extension (row 6) - fixed index of json city object buildings
$ws.Range("C6").Value = "PaduaRestricted"

# Json mode (row 9)
$ws.Range("C9").Value = "cityjson"

# District file extension (row 6)
$ws.Range("B6").Value = ".json"

# End-uses data filename (row 4)
$ws.Range("C4").Value = "ScheduleComp"

# End-uses mode (row 8)
$ws.Range("C8").Value = "Yearly"

# DHW calculation (only residential) - now a boolean FALSE (row 36)
$ws.Range("B36").Value = $false

# Selection / view state
$ws.Range("C8").Select() | Out-Null
